$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H3").Value = 79227.164
$ws.Range("J3").Value = 79227.164
$ws.Range("L3").Value = 79227.164
$ws.Range("N3").Value = -79455.164
$ws.Range("H15").Value = 1670.875
$ws.Range("I15").Value = 1670.875
$ws.Range("K15").Value = 5012.625
$ws.Range("M15").Value = -4843.625
$ws.Range("H33").Value = 231.4
$ws.Range("I33").Value = 248.85
$ws.Range("J33").Value = 161.6
$ws.Range("K33").Value = 248.85
$ws.Range("L33").Value = 161.6
$ws.Range("M33").Value = -19.84999999999999
$ws.Range("N33").Value = -619.6
$ws.Range("H86").Value = 2424.0557
$ws.Range("I86").Value = 1915.2222
$ws.Range("J86").Value = 2932.889
$ws.Range("K86").Value = 1915.2222
$ws.Range("L86").Value = 2932.889
$ws.Range("M86").Value = -792.2221999999999
$ws.Range("N86").Value = -5178.889
$ws.Range("H89").Value = 2424.0557
$ws.Range("I89").Value = 1915.2222
$ws.Range("J89").Value = 2932.889
$ws.Range("K89").Value = 9576.110999999999
$ws.Range("L89").Value = 14664.445
$ws.Range("M89").Value = -3960.110999999999
$ws.Range("N89").Value = -25896.445
$ws.Range("H98").Value = 760.48
$ws.Range("I98").Value = 763.05884
$ws.Range("J98").Value = 755
$ws.Range("K98").Value = 763.05884
$ws.Range("L98").Value = 755
$ws.Range("M98").Value = 734.94116
$ws.Range("N98").Value = -3751
$ws.Range("H102").Value = 79227.164
$ws.Range("J102").Value = 79227.164
$ws.Range("L102").Value = 79227.164
$ws.Range("N102").Value = -85717.164
$ws.Range("H122").Value = 760.48
$ws.Range("I122").Value = 763.05884
$ws.Range("J122").Value = 755
$ws.Range("K122").Value = 2289.17652
$ws.Range("L122").Value = 2265
$ws.Range("M122").Value = 160.82348
$ws.Range("N122").Value = -7165
$ws.Range("H125").Value = 1489.1666
$ws.Range("I125").Value = 1489.1666
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 13402.4994
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -10942.4994
$ws.Range("N125").ClearContents()
$ws.Range("H129").Value = 146880.14
$ws.Range("I129").Value = 336028.84
$ws.Range("K129").Value = 1008086.52
$ws.Range("M129").Value = -1003086.52
$ws.Range("H132").Value = 26680.545
$ws.Range("I132").Value = 2410.2354
$ws.Range("J132").Value = 109199.6
$ws.Range("K132").Value = 7230.706200000001
$ws.Range("L132").Value = 327598.8
$ws.Range("M132").Value = -4700.706200000001
$ws.Range("N132").Value = -332658.8
$ws.Range("H137").Value = 4884304.5
$ws.Range("I137").Value = 6458738
$ws.Range("J137").Value = 3559.6
$ws.Range("K137").Value = 19376214
$ws.Range("L137").Value = 10678.8
$ws.Range("M137").Value = -19373664
$ws.Range("N137").Value = -15778.8
$ws.Range("H138").Value = 7937.472
$ws.Range("I138").Value = 3276.8572
$ws.Range("J138").Value = 9062.448
$ws.Range("K138").Value = 9830.571599999999
$ws.Range("L138").Value = 27187.344
$ws.Range("M138").Value = -4690.571599999999
$ws.Range("N138").Value = -37467.344
$ws.Range("H140").Value = 39666.668
$ws.Range("J140").Value = 39666.668
$ws.Range("L140").Value = 39666.668
$ws.Range("N140").Value = -50026.668
$ws.Range("H141").Value = 5400.36
$ws.Range("I141").Value = 4590.6665
$ws.Range("K141").Value = 13771.9995
$ws.Range("M141").Value = -8591.999500000002

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 14499000
$ws.Range("I32").Value = 15390992
$ws.Range("K32").Value = 15390992
$ws.Range("M32").Value = -15390705
$ws.Range("H45").Value = 4996
$ws.Range("I45").Value = 3402.4
$ws.Range("J45").Value = 6766.6665
$ws.Range("K45").Value = 3402.4
$ws.Range("L45").Value = 6766.6665
$ws.Range("M45").Value = -3025.4
$ws.Range("N45").Value = -7520.6665
$ws.Range("H74").Value = 1942
$ws.Range("I74").Value = 1221.6666
$ws.Range("J74").Value = 3485.5715
$ws.Range("K74").Value = 1221.6666
$ws.Range("L74").Value = 3485.5715
$ws.Range("M74").Value = -347.6666
$ws.Range("N74").Value = -5233.5715
$ws.Range("H77").Value = 1942
$ws.Range("I77").Value = 1221.6666
$ws.Range("J77").Value = 3485.5715
$ws.Range("K77").Value = 6108.333000000001
$ws.Range("L77").Value = 17427.8575
$ws.Range("M77").Value = -1740.333000000001
$ws.Range("N77").Value = -26163.8575
$ws.Range("H132").Value = 2786.1875
$ws.Range("I132").Value = 2538.6
$ws.Range("K132").Value = 7615.799999999999
$ws.Range("M132").Value = -5085.799999999999

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 33721.848
$ws.Range("I20").Value = 2327.25
$ws.Range("J20").Value = 82021.234
$ws.Range("K20").Value = 2327.25
$ws.Range("L20").Value = 82021.234
$ws.Range("M20").Value = -2080.25
$ws.Range("N20").Value = -82515.234
$ws.Range("H99").Value = 1285.8
$ws.Range("I99").Value = 874.1667
$ws.Range("J99").Value = 2932.3333
$ws.Range("K99").Value = 874.1667
$ws.Range("L99").Value = 2932.3333
$ws.Range("M99").Value = 623.8333
$ws.Range("N99").Value = -5928.3333
$ws.Range("H134").Value = 3149.1177
$ws.Range("I134").Value = 3038.5557
$ws.Range("J134").Value = 3273.5
$ws.Range("K134").Value = 9115.667099999999
$ws.Range("L134").Value = 9820.5
$ws.Range("M134").Value = -6580.667099999999
$ws.Range("N134").Value = -14890.5

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 11510.546
$ws.Range("I16").Value = 9077.125
$ws.Range("J16").Value = 17999.666
$ws.Range("K16").Value = 9077.125
$ws.Range("L16").Value = 17999.666
$ws.Range("M16").Value = -8790.125
$ws.Range("N16").Value = -18573.666
$ws.Range("H31").Value = 2758.9517
$ws.Range("I31").Value = 1916.3334
$ws.Range("J31").Value = 4528.45
$ws.Range("K31").Value = 1916.3334
$ws.Range("L31").Value = 4528.45
$ws.Range("M31").Value = -1621.3334
$ws.Range("N31").Value = -5118.45
$ws.Range("H34").Value = 2758.9517
$ws.Range("I34").Value = 1916.3334
$ws.Range("J34").Value = 4528.45
$ws.Range("K34").Value = 1916.3334
$ws.Range("L34").Value = 4528.45
$ws.Range("M34").Value = -1714.3334
$ws.Range("N34").Value = -4932.45
$ws.Range("H58").Value = 2348.7273
$ws.Range("I58").Value = 1862.3125
$ws.Range("J58").Value = 3645.8333
$ws.Range("K58").Value = 1862.3125
$ws.Range("L58").Value = 3645.8333
$ws.Range("M58").Value = -1659.3125
$ws.Range("N58").Value = -4051.8333
$ws.Range("H86").Value = 30716.588
$ws.Range("I86").Value = 39886
$ws.Range("J86").Value = 20401
$ws.Range("K86").Value = 39886
$ws.Range("L86").Value = 20401
$ws.Range("M86").Value = -38763
$ws.Range("N86").Value = -22647
$ws.Range("H89").Value = 30716.588
$ws.Range("I89").Value = 39886
$ws.Range("J89").Value = 20401
$ws.Range("K89").Value = 199430
$ws.Range("L89").Value = 102005
$ws.Range("M89").Value = -193814
$ws.Range("N89").Value = -113237
$ws.Range("H113").Value = 11510.546
$ws.Range("I113").Value = 9077.125
$ws.Range("J113").Value = 17999.666
$ws.Range("K113").Value = 9077.125
$ws.Range("L113").Value = 17999.666
$ws.Range("M113").Value = -6907.125
$ws.Range("N113").Value = -22339.666
$ws.Range("H132").Value = 3904
$ws.Range("I132").Value = 3954.2
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 11862.6
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -9332.599999999999
$ws.Range("N132").Value = -13760
$ws.Range("H136").Value = 2348.7273
$ws.Range("I136").Value = 1862.3125
$ws.Range("K136").Value = 5586.9375
$ws.Range("L136").Value = 10937.4999
$ws.Range("M136").Value = -3036.9375
$ws.Range("N136").Value = -16037.4999

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 1545.2307
$ws.Range("I68").Value = 1462.2727
$ws.Range("J68").Value = 2001.5
$ws.Range("K68").Value = 4386.8181
$ws.Range("L68").Value = 6004.5
$ws.Range("M68").Value = -3575.8181
$ws.Range("N68").Value = -7626.5
$ws.Range("H71").Value = 1545.2307
$ws.Range("I71").Value = 1462.2727
$ws.Range("J71").Value = 2001.5
$ws.Range("K71").Value = 13160.4543
$ws.Range("L71").Value = 18013.5
$ws.Range("M71").Value = -9104.454299999999
$ws.Range("N71").Value = -26125.5
$ws.Range("H107").Value = 2527.389
$ws.Range("I107").Value = 2127.182
$ws.Range("J107").Value = 3156.2856
$ws.Range("K107").Value = 6381.545999999999
$ws.Range("L107").Value = 9468.856800000001
$ws.Range("M107").Value = -4461.545999999999
$ws.Range("N107").Value = -13308.8568
$ws.Range("H112").Value = 2400
$ws.Range("I112").Value = 2500
$ws.Range("J112").Value = 2000
$ws.Range("K112").Value = 7500
$ws.Range("L112").Value = 6000
$ws.Range("M112").Value = -6392
$ws.Range("N112").Value = -8216
$ws.Range("H115").Value = 1875
$ws.Range("I115").Value = 500
$ws.Range("J115").Value = 3250
$ws.Range("K115").Value = 1500
$ws.Range("L115").Value = 9750
$ws.Range("M115").Value = -325
$ws.Range("N115").Value = -12100
$ws.Range("H137").Value = 1633
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H99").Value = 22400
$ws.Range("I99").Value = 11666.667
$ws.Range("J99").Value = 38500
$ws.Range("K99").Value = 11666.667
$ws.Range("L99").Value = 38500
$ws.Range("M99").Value = -9420.666999999999
$ws.Range("N99").Value = -42992
$ws.Range("H113").Value = 17111.111
$ws.Range("I113").Value = 9800
$ws.Range("J113").Value = 26250
$ws.Range("K113").Value = 9800
$ws.Range("L113").Value = 26250
$ws.Range("M113").Value = -7630
$ws.Range("N113").Value = -30590
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value = 1808.6666
$ws.Range("I113").Value = 750.6667
$ws.Range("J113").Value = 2866.6667
$ws.Range("K113").Value = 2252.0001
$ws.Range("L113").Value = 8600.000100000001
$ws.Range("M113").Value = -82.0001000000002
$ws.Range("N113").Value = -12940.0001
$ws.Range("H132").Value = 3364.32
$ws.Range("I132").Value = 3054.238
$ws.Range("K132").Value = 9162.714
$ws.Range("M132").Value = -6632.714
$ws.Range("H135").Value = 46888.89
$ws.Range("J135").Value = 46888.89
$ws.Range("L135").Value = 46888.89
$ws.Range("N135").Value = -57028.89
$ws.Range("H141").Value = 85177.37
$ws.Range("J141").Value = 87630.10000000001
$ws.Range("L141").Value = 87630.10000000001
$ws.Range("N141").Value = -97990.10000000001
